$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 corresponds to REVOLSL_bil_usd (Revolving Consumer Credit) - auto-updated data + news
$ws.Range("E6").Value = 1328559.61
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "Dec 2025"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = 1101106.887024793
$ws.Range("H6").Value = 31594.05000000005
$ws.Range("I6").Value = 0.02435997606597977
